# Updated symbol list values (Price / Volume(1h)) per the diff.
# Values are entered with a leading apostrophe (quote-prefix), exactly as Excel
# keeps numeric-looking text literal when a user types an apostrophe-prefixed
# entry -- this preserves the original inline-string cell type/value without
# forcing a different NumberFormat on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.95"
$ws.Range("E2").Value = "'-0.44%"
$ws.Range("D3").Value = "'41.97"
$ws.Range("E3").Value = "'7.51%"
$ws.Range("D4").Value = "'5.666"
$ws.Range("E4").Value = "'-1.54%"
$ws.Range("D5").Value = "'0.08342"
$ws.Range("E5").Value = "'3.76%"
$ws.Range("D6").Value = "'2.037"
$ws.Range("E6").Value = "'3.92%"
$ws.Range("D7").Value = "'8.800"
$ws.Range("E7").Value = "'1.90%"
$ws.Range("E8").Value = "'0.43%"
$ws.Range("D10").Value = "'0.9311"
$ws.Range("D11").Value = "'0.1292"
$ws.Range("E11").Value = "'2.10%"
$ws.Range("D12").Value = "'0.1966"
$ws.Range("E12").Value = "'0.96%"
$ws.Range("D13").Value = "'0.09368"
$ws.Range("E13").Value = "'1.54%"
$ws.Range("D14").Value = "'0.03916"
$ws.Range("E14").Value = "'9.50%"
$ws.Range("D15").Value = "'0.1061"
$ws.Range("E15").Value = "'0.64%"
$ws.Range("E16").Value = "'-0.29%"
$ws.Range("D17").Value = "'0.006182"
$ws.Range("E17").Value = "'-1.42%"
$ws.Range("D18").Value = "'3.447"
$ws.Range("E18").Value = "'2.25%"
$ws.Range("E19").Value = "'1.50%"
$ws.Range("D20").Value = "'8.328"
$ws.Range("E20").Value = "'-4.79%"
$ws.Range("D21").Value = "'0.1360"
$ws.Range("E21").Value = "'1.39%"
$ws.Range("D22").Value = "'0.2455"
$ws.Range("E22").Value = "'-8.42%"
$ws.Range("D23").Value = "'0.04414"
$ws.Range("E23").Value = "'-0.73%"
$ws.Range("D24").Value = "'0.001249"
$ws.Range("E24").Value = "'-1.47%"
$ws.Range("D25").Value = "'0.004378"
$ws.Range("E25").Value = "'-1.77%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.68%"
$ws.Range("D39").Value = "'0.02833"
$ws.Range("E39").Value = "'12.42%"
$ws.Range("D40").Value = "'0.05542"
$ws.Range("E40").Value = "'1.27%"
$ws.Range("D41").Value = "'0.007809"
$ws.Range("E41").Value = "'4.76%"
$ws.Range("D42").Value = "'0.1441"
$ws.Range("E42").Value = "'2.38%"
$ws.Range("D43").Value = "'0.008941"
$ws.Range("E43").Value = "'-10.09%"
$ws.Range("D44").Value = "'0.002241"
$ws.Range("E44").Value = "'11.38%"
$ws.Range("D45").Value = "'0.01174"
$ws.Range("D46").Value = "'0.00007016"
$ws.Range("E46").Value = "'3.15%"
$ws.Range("E47").Value = "'-0.68%"
$ws.Range("D48").Value = "'0.003176"
$ws.Range("E48").Value = "'4.33%"
$ws.Range("E49").Value = "'-0.50%"
$ws.Range("E50").Value = "'-0.68%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.68%"
